# Updated for new routine
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: replace the header row with a single title cell ---
$ws.Range("A1:L1").ClearContents()
$ws.Range("A1").Value = "Herald College Kathmandu"

# --- Data rows 2-10: reorder/rewrite columns C:J, drop K:L (Level/Course) ---
# New column layout: A Day | B Time | C Hours | D Module Code | E Module Title |
#                     F Class Type | G Lecturer | H Group | I Block | J Room

$data = @(
    @{Row=2;  Hours=2.5; Code="5CS024"; Title="Collaborative Development";                 Group="L5CG5";          Block="HCK"; Room="Lab-04 Patan"},
    @{Row=3;  Hours=2;   Code="5CS022"; Title="Human Computer Interaction";                 Group="L5CG(5+6+7+8)";  Block="WLV"; Room="LT-02 Telford"},
    @{Row=4;  Hours=2;   Code="5CS020"; Title="Distributed and Cloud Systems Programming";  Group="L5CG(5+6+7+8)";  Block="WLV"; Room="LT-01 Wulfruna"},
    @{Row=5;  Hours=2;   Code="5CS024"; Title="Collaborative Development";                  Group="L5CG(5+6+7+8)";  Block="WLV"; Room="LT-02 Telford"},
    @{Row=6;  Hours=2;   Code="5CS020"; Title="Distributed and Cloud Systems Programming";  Group="L5CG5";          Block="WLV"; Room="TR-01 Dudley"},
    @{Row=7;  Hours=2;   Code="5CS022"; Title="Human Computer Interaction";                 Group="L5CG5";          Block="WLV"; Room="TR-01 Dudley"},
    @{Row=8;  Hours=2.5; Code="5CS020"; Title="Distributed and Cloud Systems Programming";  Group="L5CG5";          Block="WLV"; Room="Lab-02 Moseley"},
    @{Row=9;  Hours=2;   Code="5CS024"; Title="Collaborative Development";                  Group="L5CG5";          Block="WLV"; Room="TR-03 Westbromwich"},
    @{Row=10; Hours=2.5; Code="5CS022"; Title="Human Computer Interaction";                 Group="L5CG5";          Block="HCK"; Room="TR-11 Nagarjung"}
)

foreach ($d in $data) {
    $r = $d.Row
    $ws.Cells.Item($r, 3).Value = $d.Hours
    $ws.Cells.Item($r, 4).Value = $d.Code
    $ws.Cells.Item($r, 5).Value = $d.Title
    $ws.Cells.Item($r, 8).Value = $d.Group
    $ws.Cells.Item($r, 9).Value = $d.Block
    $ws.Cells.Item($r, 10).Value = $d.Room
    $ws.Range("K$r`:L$r").ClearContents()
}

$ws.Range("K1:L10").ClearContents()
